$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add the two new "authorization window" columns: StartTime / EndTime
# ---------------------------------------------------------------------
$ws.Range("J1").Value = "StartTime"
$ws.Range("K1").Value = "EndTime"

# Build the date/time style once on a template cell (value + number
# format + font size) so the engine records a single combined cell
# style, then replicate that style (format only) onto the remaining
# cells with Copy / PasteSpecial before filling in their own values.
$template = $ws.Range("J2")
$template.Value = 44018.999988425923
$template.NumberFormat = "yyyy/mm/dd\ h:mm:ss"
$template.Font.Size = 10.5

$template.Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("K2").Value = 44171.999988425923
$ws.Range("J3").Value = 44018.999988425923
$ws.Range("K3").Value = 44171.999988425923

# ---------------------------------------------------------------------
# Column widths for the two new columns
# ---------------------------------------------------------------------
$ws.Columns("J").ColumnWidth = 18.666666666666668
$ws.Columns("K").ColumnWidth = 17.166666666666668

# ---------------------------------------------------------------------
# Selection ends up on K4 after entering the data
# ---------------------------------------------------------------------
$ws.Range("K4").Select() | Out-Null
